{"js": "// Update the two-digit \u00f7 one-digit division answers shown in the table.\n// Each entry identifies the target cell by its (row, col) position in the\n// table (rows include the blank spacer rows that separate the answer rows)\n// plus the exact old/new text, so that cells sharing identical text\n// elsewhere in the table are never touched by mistake.\nconst changes = [\n  { row: 0, col: 0, oldText: \"22\u00f75=4, 2\", newText: \"39\u00f75=7, 4\" },\n  { row: 0, col: 1, oldText: \"96\u00f76=16, 0\", newText: \"72\u00f75=14, 2\" },\n  { row: 0, col: 2, oldText: \"50\u00f74=12, 2\", newText: \"26\u00f79=2, 8\" },\n  { row: 0, col: 3, oldText: \"99\u00f79=11, 0\", newText: \"98\u00f76=16, 2\" },\n  { row: 0, col: 4, oldText: \"50\u00f78=6, 2\", newText: \"46\u00f78=5, 6\" },\n  { row: 4, col: 0, oldText: \"26\u00f75=5, 1\", newText: \"13\u00f77=1, 6\" },\n  { row: 4, col: 1, oldText: \"81\u00f77=11, 4\", newText: \"67\u00f74=16, 3\" },\n  { row: 4, col: 2, oldText: \"75\u00f75=15, 0\", newText: \"16\u00f74=4, 0\" },\n  { row: 4, col: 3, oldText: \"35\u00f77=5, 0\", newText: \"51\u00f77=7, 2\" },\n  { row: 4, col: 4, oldText: \"13\u00f76=2, 1\", newText: \"20\u00f77=2, 6\" },\n  { row: 8, col: 0, oldText: \"66\u00f74=16, 2\", newText: \"15\u00f76=2, 3\" },\n  { row: 8, col: 1, oldText: \"73\u00f78=9, 1\", newText: \"73\u00f73=24, 1\" },\n  { row: 8, col: 3, oldText: \"15\u00f79=1, 6\", newText: \"94\u00f77=13, 3\" },\n  { row: 8, col: 4, oldText: \"84\u00f74=21, 0\", newText: \"30\u00f72=15, 0\" },\n  { row: 12, col: 0, oldText: \"17\u00f78=2, 1\", newText: \"23\u00f78=2, 7\" },\n  { row: 12, col: 1, oldText: \"15\u00f73=5, 0\", newText: \"78\u00f76=13, 0\" },\n  { row: 12, col: 2, oldText: \"94\u00f74=23, 2\", newText: \"43\u00f73=14, 1\" },\n  { row: 12, col: 3, oldText: \"54\u00f76=9, 0\", newText: \"65\u00f78=8, 1\" },\n  { row: 12, col: 4, oldText: \"16\u00f74=4, 0\", newText: \"64\u00f79=7, 1\" },\n  { row: 16, col: 0, oldText: \"37\u00f74=9, 1\", newText: \"24\u00f75=4, 4\" },\n  { row: 16, col: 1, oldText: \"66\u00f74=16, 2\", newText: \"14\u00f73=4, 2\" },\n  { row: 16, col: 2, oldText: \"41\u00f73=13, 2\", newText: \"75\u00f77=10, 5\" },\n  { row: 16, col: 3, oldText: \"10\u00f79=1, 1\", newText: \"47\u00f76=7, 5\" },\n  { row: 16, col: 4, oldText: \"48\u00f78=6, 0\", newText: \"36\u00f74=9, 0\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const change of changes) {\n  const cell = table.getCell(change.row, change.col);\n  // Scope the search to this single cell's body so that duplicate answer\n  // strings elsewhere in the table can't be matched/replaced by accident.\n  const results = cell.body.search(change.oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(change.newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();", "ps1": "# Update the two-digit \u00f7 one-digit division answers shown in the table.\n# The table has 20 rows: one answer row followed by 3 blank spacer rows,\n# repeated 5 times (answer rows are COM rows 1, 5, 9, 13, 17).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Replace-CellText($row, $col, $oldText, $newText) {\n  $cell = $tbl.Cell($row, $col)\n  $cellRange = $cell.Range\n  # Re-wrap the cell's Start/End as a fresh document Range: Find on the\n  # Cell's own Range object searches from the top of the document in this\n  # host, while a Range built from $d.Range(start, end) is correctly\n  # confined to the cell - this keeps duplicate answer strings elsewhere\n  # in the table from being matched/replaced by accident.\n  $rng = $d.Range($cellRange.Start, $cellRange.End)\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  # wdFindStop=0, wdReplaceOne=1\n  $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 1)\n  if (-not $result) {\n    Write-Output \"WARNING: replace failed for cell ($row,$col): '$oldText' -> '$newText'\"\n  }\n}\n\nReplace-CellText 1 1 \"22\u00f75=4, 2\" \"39\u00f75=7, 4\"\nReplace-CellText 1 2 \"96\u00f76=16, 0\" \"72\u00f75=14, 2\"\nReplace-CellText 1 3 \"50\u00f74=12, 2\" \"26\u00f79=2, 8\"\nReplace-CellText 1 4 \"99\u00f79=11, 0\" \"98\u00f76=16, 2\"\nReplace-CellText 1 5 \"50\u00f78=6, 2\" \"46\u00f78=5, 6\"\nReplace-CellText 5 1 \"26\u00f75=5, 1\" \"13\u00f77=1, 6\"\nReplace-CellText 5 2 \"81\u00f77=11, 4\" \"67\u00f74=16, 3\"\nReplace-CellText 5 3 \"75\u00f75=15, 0\" \"16\u00f74=4, 0\"\nReplace-CellText 5 4 \"35\u00f77=5, 0\" \"51\u00f77=7, 2\"\nReplace-CellText 5 5 \"13\u00f76=2, 1\" \"20\u00f77=2, 6\"\nReplace-CellText 9 1 \"66\u00f74=16, 2\" \"15\u00f76=2, 3\"\nReplace-CellText 9 2 \"73\u00f78=9, 1\" \"73\u00f73=24, 1\"\nReplace-CellText 9 4 \"15\u00f79=1, 6\" \"94\u00f77=13, 3\"\nReplace-CellText 9 5 \"84\u00f74=21, 0\" \"30\u00f72=15, 0\"\nReplace-CellText 13 1 \"17\u00f78=2, 1\" \"23\u00f78=2, 7\"\nReplace-CellText 13 2 \"15\u00f73=5, 0\" \"78\u00f76=13, 0\"\nReplace-CellText 13 3 \"94\u00f74=23, 2\" \"43\u00f73=14, 1\"\nReplace-CellText 13 4 \"54\u00f76=9, 0\" \"65\u00f78=8, 1\"\nReplace-CellText 13 5 \"16\u00f74=4, 0\" \"64\u00f79=7, 1\"\nReplace-CellText 17 1 \"37\u00f74=9, 1\" \"24\u00f75=4, 4\"\nReplace-CellText 17 2 \"66\u00f74=16, 2\" \"14\u00f73=4, 2\"\nReplace-CellText 17 3 \"41\u00f73=13, 2\" \"75\u00f77=10, 5\"\nReplace-CellText 17 4 \"10\u00f79=1, 1\" \"47\u00f76=7, 5\"\nReplace-CellText 17 5 \"48\u00f78=6, 0\" \"36\u00f74=9, 0\"\n"}
